# Historical.xlsx update: add "Meet Basic Needs" and "Open Ended Questions"
# domain rows to the question bank, plus trailing blank rows, matching the
# author's "draft historical and updated social support" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row data -----------------------------------------------------------
# Each tuple: (row number, Question text, Answers text, Domain text, row height)
$domainBasicNeeds = "Meet Basic Needs for Health and Well-Being"
$domainOpenEnded  = "Open Ended Questions"
$openAnswer       = "Open answer"

$rows = @(
  @{ Row=3;  Q="Prior to coronavirus (COVID-19), did you get free or reduced lunches for your child(ren)? ";
     A="• Yes, my child(ren) received free or reduced lunches  `n• No, my child did not receive free or reduced lunches, but they were available `n• No, free or reduced lunches were not available for my children  `n• No, but I had planned on getting free or reduced lunches for my children soon  `n• No, I didn't know how to access this resource `n• Not applicable ";
     D=$domainBasicNeeds; H=409.6 }
  @{ Row=4;  Q="Do you currently get free or reduced lunches for your child(ren)? ";
     A="• Yes, my child(ren) receive(s) free or reduced lunches  `n• No, my child does not receive free or reduced lunches, but they were available `n• No, free or reduced lunches are not available for my children  `n• No, but I plan on getting free or reduced lunches for my children soon  `n• No, I don't know how to access this resource `n• Not applicable ";
     D=$domainBasicNeeds; H=409.6 }
  @{ Row=5;  Q="In the past month, how hard has it been for you to pay for the very basics like food, housing, medical care, and heating? ";
     A="• Very hard `n• Hard `n• Somewhat hard  `n• Not very hard ";
     D=$domainBasicNeeds; H=204 }
  @{ Row=6;  Q="Which of these needs have been hard to pay for in the past month? Select all that apply.  ";
     A="• Food `n• Housing `n• Utilities (electric, water, trash) `n• Healthcare `n• Social `n• Emotional `n• Childcare `n• Other (please specify) [text entry] `n• None of the above";
     D=$domainBasicNeeds; H=306 }
  @{ Row=7;  Q="What are the biggest challenges and concerns for you and your family right now? ";
     A=$openAnswer; D=$domainOpenEnded; H=136 }
  @{ Row=8;  Q="What is helping you and your family the most right now? ";
     A=$openAnswer; D=$domainOpenEnded; H=102 }
  @{ Row=9;  Q="What is on your mind the most when you think about your community re-opening? ";
     A=$openAnswer; D=$domainOpenEnded; H=153 }
  @{ Row=10; Q="What concerns do you have about your place of employment and/or your child’s child care setting re-opening? ";
     A=$openAnswer; D=$domainOpenEnded; H=187 }
  @{ Row=11; Q="Is there anything else you would like to tell us about you and your family’s experiences during the COVID-19 pandemic? ";
     A=$openAnswer; D=$domainOpenEnded; H=204 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    # Copy the formatting from row 2 (A2:C2) down onto the new row first, so
    # the new cells pick up the same wrap/alignment styles (s="1"/s="1"/s="2")
    # without registering extra unused cell styles.
    $ws.Range("A2:C2").Copy()
    $ws.Range("A$rowNum`:C$rowNum").PasteSpecial(-4122)

    $ws.Range("A$rowNum").Value = $r.Q
    $ws.Range("B$rowNum").Value = $r.A
    $ws.Range("C$rowNum").Value = $r.D

    $ws.Rows.Item($rowNum).RowHeight = $r.H
}

# --- Trailing blank rows (12-17), same formatting, no values -------------
for ($rowNum = 12; $rowNum -le 17; $rowNum++) {
    $ws.Range("A2:C2").Copy()
    $ws.Range("A$rowNum`:C$rowNum").PasteSpecial(-4122)
}

$excel.CutCopyMode = 0

# --- Selection: rows 3:17 selected, matching the author's saved view -----
$ws.Range("A3:A17").EntireRow.Select() | Out-Null
